$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffix = "  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on. "

# Row 10 already has the MCA comment ("Nov 2024 data...") and is left untouched.
# Append the MCA comment to every other dataset's "LatestPeriod" (column B) cell, rows 2-24.
for ($r = 2; $r -le 24; $r++) {
    if ($r -eq 10) {
        continue
    }
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value2
    if ($current -eq $null) {
        continue
    }
    $cell.Value2 = $current + "." + $suffix
}
